$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 6 with the new time-tracking entry
$ws.Range("A6").Value = "Tuesday 17.4.18"
$ws.Range("B6").Value = "1030 - 1130"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "Implementing LinkedList"

# Update the active selection to D6, matching the saved cursor position
$ws.Range("D6").Select()

# Force recalculation so the SUBTOTAL formula in C32 reflects the new row
$excel.CalculateFullRebuild()
